$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Merge the "We can make..." paragraph with the "And beforehand..."
#    paragraph (which loses its own pPr) and rewrite the sentence.
# ------------------------------------------------------------------
$pA = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text.StartsWith("We can make the group of employees")) {
        $pA = $cand
        break
    }
}
$pAIndex = $pA.Index
$pB = $d.Paragraphs.Item($pAIndex + 1)
$pBTextRange = $d.Range($pB.Range.Start, $pB.Range.End - 1)
$pBText = $pBTextRange.Text

$pAAgain = $d.Paragraphs.Item($pAIndex)
$insertPos = $pAAgain.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter($pBText)

$pBFresh = $d.Paragraphs.Item($pAIndex + 1)
$pBFresh.Range.Delete()

$mergedPara = $d.Paragraphs.Item($pAIndex)
$mergedRange = $mergedPara.Range
$mergedRange.Find.Execute("We can make the group of employees with respective to their Resignation month.And beforehand we will transfer or hire the new employees for the work to be assigned.", $false, $false, $false, $false, $false, $true, 1, $false, "With the data which collected from client, we can predict the next person to resign.", 2)

# ------------------------------------------------------------------
# 2) "Machine Learning " Unsupervised -> " supervised-Classification"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Machine Learning " + [char]0x2013 + " Unsupervised", $false, $false, $false, $false, $false, $true, 1, $false, "Machine Learning " + [char]0x2013 + " supervised-Classification", 2)

# ------------------------------------------------------------------
# 3) "Sample Dataset which has only I/Ps." -> "Sample Dataset."
# ------------------------------------------------------------------
$d.Content.Find.Execute("Sample Dataset which has only I/Ps.", $false, $false, $false, $false, $false, $true, 1, $false, "Sample Dataset.", 2)

# ------------------------------------------------------------------
# 4) Table: add a 3rd "Output" column, edit header/cell text.
# ------------------------------------------------------------------
$t = $d.Tables.Item(1)
$t.Columns.Add() | Out-Null
$t.Columns.Item(3).Width = 144

# Header row text swaps.
$d.Content.Find.Execute("Input 1", $false, $false, $false, $false, $false, $true, 1, $false, "Emp.Name", 2)
$d.Content.Find.Execute("Input 2(resignation date)", $false, $false, $false, $false, $false, $true, 1, $false, "Age", 2)
$d.Content.Find.Execute("March 31", $false, $false, $false, $false, $false, $true, 1, $false, "32", 2)
$d.Content.Find.Execute("April 13", $false, $false, $false, $false, $false, $true, 1, $false, "27", 2)
$d.Content.Find.Execute("March 12", $false, $false, $false, $false, $false, $true, 1, $false, "29", 2)

# New "Output" column cells.
$cell1 = $t.Cell(1, 3)
$cell1.Range.Text = "Output"
$cell1.Range.Font.Bold = $true
$cell1.Range.Font.Size = 11
$cell1.Range.Font.SizeBi = 11

$cell2 = $t.Cell(2, 3)
$cell2.Range.Text = "Resign"
$cell2.Range.Font.Bold = $true
$cell2.Range.Font.Size = 11
$cell2.Range.Font.SizeBi = 11

$cell3 = $t.Cell(3, 3)
$cell3.Range.Text = "Working"
$cell3.Range.Font.Bold = $true
$cell3.Range.Font.Size = 11
$cell3.Range.Font.SizeBi = 11

$cell4 = $t.Cell(4, 3)
$cell4.Range.Text = "working"
$cell4.Range.Font.Bold = $true
$cell4.Range.Font.Size = 11
$cell4.Range.Font.SizeBi = 11

Write-Output "done"
